# Update automàtic: dades i banners [2026-02-16 18:50]
# Refresh scraped MeteoCat rows: extraction timestamps + updated readings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-16 18:48:25"
$ws.Range("I2").Value = "20.0 mm"
$ws.Range("E3").Value = "2026-02-16 18:48:27"
$ws.Range("I3").Value = "9.9 mm"
$ws.Range("E4").Value = "2026-02-16 18:48:30"
$ws.Range("E5").Value = "2026-02-16 18:48:32"
$ws.Range("I5").Value = "23.5 mm"
$ws.Range("N5").Value = "-1.6 °C 18:27 TU"
$ws.Range("E6").Value = "2026-02-16 18:48:35"
$ws.Range("E7").Value = "2026-02-16 18:48:37"
$ws.Range("E8").Value = "2026-02-16 18:48:40"
$ws.Range("E9").Value = "2026-02-16 18:48:42"
$ws.Range("H9").Value = "'72%"
$ws.Range("O9").Value = "11.2 °C"
$ws.Range("E10").Value = "2026-02-16 18:48:44"
$ws.Range("E11").Value = "2026-02-16 18:48:45"
$ws.Range("H11").Value = "'78%"
$ws.Range("O11").Value = "6.7 °C"
$ws.Range("E12").Value = "2026-02-16 18:48:47"
$ws.Range("H12").Value = "'81%"
$ws.Range("O12").Value = "10.7 °C"
$ws.Range("E13").Value = "2026-02-16 18:48:48"
$ws.Range("E14").Value = "2026-02-16 18:48:49"
$ws.Range("O14").Value = "16.1 °C"
$ws.Range("E15").Value = "2026-02-16 18:48:50"
$ws.Range("O15").Value = "11.3 °C"
$ws.Range("E16").Value = "2026-02-16 18:48:51"
$ws.Range("L16").Value = "102.6 km/h - 206º 18:26 TU"
$ws.Range("N16").Value = "-1.6 °C 18:23 TU"
$ws.Range("E17").Value = "2026-02-16 18:48:52"
$ws.Range("N17").Value = "3.8 °C 18:10 TU"
$ws.Range("O17").Value = "6.1 °C"
$ws.Range("E18").Value = "2026-02-16 18:48:53"
$ws.Range("H18").Value = "'73%"
$ws.Range("O18").Value = "11.0 °C"
$ws.Range("E19").Value = "2026-02-16 18:48:54"
$ws.Range("E20").Value = "2026-02-16 18:48:55"
$ws.Range("I20").Value = "0.3 mm"
$ws.Range("E21").Value = "2026-02-16 18:48:57"
$ws.Range("L21").Value = "50.4 km/h - 318º 18:20 TU"
$ws.Range("O21").Value = "8.4 °C"
$ws.Range("E22").Value = "2026-02-16 18:48:59"
$ws.Range("E23").Value = "2026-02-16 18:49:01"
$ws.Range("I23").Value = "13.2 mm"
$ws.Range("L23").Value = "79.6 km/h - 266º 18:27 TU"
$ws.Range("E24").Value = "2026-02-16 18:49:04"
$ws.Range("E25").Value = "2026-02-16 18:49:06"
$ws.Range("I25").Value = "5.6 mm"
$ws.Range("E26").Value = "2026-02-16 18:49:09"
$ws.Range("E27").Value = "2026-02-16 18:49:11"
$ws.Range("E28").Value = "2026-02-16 18:49:13"
$ws.Range("H28").Value = "'71%"
$ws.Range("O28").Value = "9.6 °C"
$ws.Range("E29").Value = "2026-02-16 18:49:15"
$ws.Range("E30").Value = "2026-02-16 18:49:18"
$ws.Range("E31").Value = "2026-02-16 18:49:20"
$ws.Range("E32").Value = "2026-02-16 18:49:22"
$ws.Range("E33").Value = "2026-02-16 18:49:25"
$ws.Range("I33").Value = "0.7 mm"
$ws.Range("J33").Value = "1013.9 hPa"
$ws.Range("O33").Value = "6.3 °C"
$ws.Range("E34").Value = "2026-02-16 18:49:27"
$ws.Range("N34").Value = "2.2 °C 18:29 TU"
$ws.Range("E35").Value = "2026-02-16 18:49:30"
$ws.Range("J35").Value = "1016.8 hPa"
$ws.Range("E36").Value = "2026-02-16 18:49:32"
$ws.Range("E37").Value = "2026-02-16 18:49:35"
$ws.Range("J37").Value = "1014.8 hPa"
$ws.Range("O37").Value = "6.6 °C"
$ws.Range("E38").Value = "2026-02-16 18:49:37"
$ws.Range("O38").Value = "12.0 °C"
$ws.Range("E39").Value = "2026-02-16 18:49:40"
$ws.Range("I39").Value = "3.2 mm"
$ws.Range("N39").Value = "-1.3 °C 18:28 TU"
$ws.Range("E40").Value = "2026-02-16 18:49:42"
$ws.Range("J40").Value = "1016.5 hPa"
$ws.Range("O40").Value = "6.9 °C"
$ws.Range("E41").Value = "2026-02-16 18:49:45"
$ws.Range("J41").Value = "1014.9 hPa"
$ws.Range("O41").Value = "17.5 °C"
$ws.Range("E42").Value = "2026-02-16 18:49:47"
$ws.Range("O42").Value = "11.3 °C"
$ws.Range("E43").Value = "2026-02-16 18:49:49"
$ws.Range("O43").Value = "8.4 °C"
$ws.Range("E44").Value = "2026-02-16 18:49:52"
$ws.Range("E45").Value = "2026-02-16 18:49:54"
$ws.Range("I45").Value = "16.7 mm"
$ws.Range("E46").Value = "2026-02-16 18:49:57"
$ws.Range("J46").Value = "1017.0 hPa"
$ws.Range("K46").Value = "12.7 MJ/m2"
